$d = $word.ActiveDocument

# --- 1) Typo fix: "EMG" -> split into "EM" + "G" runs, with the
#        "_GoBack" bookmark (auto-inserted by Word at the last edit
#        point) landing between them. Adding a bookmark with the name
#        "_GoBack" also removes any pre-existing "_GoBack" bookmark
#        elsewhere in the document (Word only keeps one), which is
#        exactly the behaviour captured by the diff.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("record surface EMG (", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$hit = $find.Parent.Duplicate
# $hit now spans "record surface EMG (" ; locate "EMG" inside it.
$emgStart = $hit.Start + "record surface ".Length
$emgEnd = $emgStart + "EMG".Length
$splitPoint = $emgStart + "EM".Length

$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# --- 2) Merge " we discuss " + "the design and control of AWS, which "
#        into a single run by replacing the combined text with itself.
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Execute(" we discuss the design and control of AWS, which ", $true, $false, $false, $false, $false, $true, 1, $false, " we discuss the design and control of AWS, which ", 2) | Out-Null

# --- 3) Merge " on Unplugged Powered Suit (UPS)" + " by keeping human "
#        into a single run the same way.
$find3 = $d.Content.Find
$find3.ClearFormatting()
$find3.Replacement.ClearFormatting()
$find3.Execute(" on Unplugged Powered Suit (UPS) by keeping human ", $true, $false, $false, $false, $false, $true, 1, $false, " on Unplugged Powered Suit (UPS) by keeping human ", 2) | Out-Null
